$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new row at 29 for the new item "COLD FREE 20 TAB."
#    (alphabetically between "COLCHICINE 500MCG 100 TAB" and
#    "CONTAFEVER N 200MG/5ML SUSP. 120ML"). This shifts the old
#    rows 29-129 down to 30-130.
$ws.Rows("29:29").Insert()

# 2) Copy formatting (styles) from the row just below (old row 29,
#    now at row 30) onto the freshly inserted, blank row 29.
$ws.Range("A30:N30").Copy()
$ws.Range("A29:N29").PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0

# Keep the same row height the table already uses for this band.
$ws.Rows("29:29").RowHeight = 24.75

# 3) Re-create the merged cells for the new row (B:G, H:K, L:M),
#    matching every other data row in the table.
$ws.Range("B29:G29").Merge()
$ws.Range("H29:K29").Merge()
$ws.Range("L29:M29").Merge()

# 4) Populate the new row's values.
$ws.Range("A29").Value = 26
$ws.Range("B29").Value = "COLD FREE 20 TAB."
$ws.Range("H29").Value = "1:0"
$ws.Range("L29").Value = 25
$ws.Range("N29").Value = "0:2"

# 5) The "م" (sequence number) column is a static count, not a
#    formula, so every row pushed down by the insert needs its
#    number bumped by one to stay sequential (26..125).
for ($r = 30; $r -le 128; $r++) {
    $cur = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($r, 1).Value = $cur + 1
}

# 6) The grand-total cell (column K) is a literal sum, not a
#    formula either, so add the new item's price (25) to it.
#    It used to live at K128 and is now at K129 after the insert.
$oldTotal = $ws.Range("K129").Value2
$ws.Range("K129").Value = $oldTotal + 25
